# Auto-generated edit script: updates market-price derived
# columns (H-N) on several worksheets, applying refreshed
# currentAveragePrice-based figures pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6
$ws.Cells.Item(6, 8).Value = 2164.9375
$ws.Cells.Item(6, 9).Value = 91.28570999999999
$ws.Cells.Item(6, 11).Value = 273.85713
$ws.Cells.Item(6, 13).Value = -161.85713

# row 107
$ws.Cells.Item(107, 8).Value = 1385
$ws.Cells.Item(107, 9).Value = 1699.3334
$ws.Cells.Item(107, 10).Value = 1127.8182
$ws.Cells.Item(107, 11).Value = 1699.3334
$ws.Cells.Item(107, 12).Value = 1127.8182
$ws.Cells.Item(107, 13).Value = 220.6666
$ws.Cells.Item(107, 14).Value = -4967.8182

# row 137
$ws.Cells.Item(137, 8).Value = 1444951.4
$ws.Cells.Item(137, 9).Value = 2072004
$ws.Cells.Item(137, 10).Value = 2730.2
$ws.Cells.Item(137, 11).Value = 6216012
$ws.Cells.Item(137, 12).Value = 8190.599999999999
$ws.Cells.Item(137, 13).Value = -6213462
$ws.Cells.Item(137, 14).Value = -13290.6

$ws = $wb.Worksheets.Item("ARM")
# row 3
$ws.Cells.Item(3, 8).Value = 11596.2
$ws.Cells.Item(3, 9).Value = 10000
$ws.Cells.Item(3, 10).Value = 11995.25
$ws.Cells.Item(3, 11).Value = 10000
$ws.Cells.Item(3, 12).Value = 11995.25
$ws.Cells.Item(3, 13).Value = -9885
$ws.Cells.Item(3, 14).Value = -12225.25

# row 32
$ws.Cells.Item(32, 8).Value = 5640.86
$ws.Cells.Item(32, 9).Value = 4255.405
$ws.Cells.Item(32, 11).Value = 4255.405
$ws.Cells.Item(32, 13).Value = -3968.405

# row 48
$ws.Cells.Item(48, 8).Value = 79800
$ws.Cells.Item(48, 10).Value = 79800
$ws.Cells.Item(48, 12).Value = 79800
$ws.Cells.Item(48, 14).Value = -80568

# row 137
$ws.Cells.Item(137, 8).Value = 39191.25
$ws.Cells.Item(137, 10).Value = 40532.855
$ws.Cells.Item(137, 12).Value = 40532.855
$ws.Cells.Item(137, 14).Value = -50732.855

$ws = $wb.Worksheets.Item("BSM")
# row 5
$ws.Cells.Item(5, 8).Value = 2351.5715
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(5, 10).Value = 4112.25
$ws.Cells.Item(5, 11).Value = 4
$ws.Cells.Item(5, 12).Value = 4112.25
$ws.Cells.Item(5, 13).Value = 109
$ws.Cells.Item(5, 14).Value = -4338.25

# row 137
$ws.Cells.Item(137, 8).Value = 58865
$ws.Cells.Item(137, 10).Value = 58865
$ws.Cells.Item(137, 12).Value = 58865
$ws.Cells.Item(137, 14).Value = -69065

$ws = $wb.Worksheets.Item("CRP")
# row 2
$ws.Cells.Item(2, 8).Value = 3001
$ws.Cells.Item(2, 9).Value = 1502
$ws.Cells.Item(2, 10).Value = 4500
$ws.Cells.Item(2, 11).Value = 1502
$ws.Cells.Item(2, 12).Value = 4500
$ws.Cells.Item(2, 13).Value = -1389
$ws.Cells.Item(2, 14).Value = -4726

# row 17
$ws.Cells.Item(17, 8).Value = 17999
$ws.Cells.Item(17, 10).Value = 17999
$ws.Cells.Item(17, 12).Value = 17999
$ws.Cells.Item(17, 14).Value = -18347

# row 25
$ws.Cells.Item(25, 8).Value = 34900
$ws.Cells.Item(25, 10).Value = 34900
$ws.Cells.Item(25, 12).Value = 34900
$ws.Cells.Item(25, 14).Value = -35248

# row 31
$ws.Cells.Item(31, 8).Value = 2551.9033
$ws.Cells.Item(31, 9).Value = 1017.0417
$ws.Cells.Item(31, 11).Value = 1017.0417
$ws.Cells.Item(31, 13).Value = -722.0417

# row 34
$ws.Cells.Item(34, 8).Value = 2551.9033
$ws.Cells.Item(34, 9).Value = 1017.0417
$ws.Cells.Item(34, 11).Value = 1017.0417
$ws.Cells.Item(34, 13).Value = -815.0417

# row 68
$ws.Cells.Item(68, 8).Value = 50167.668
$ws.Cells.Item(68, 10).Value = 50167.668
$ws.Cells.Item(68, 12).Value = 50167.668
$ws.Cells.Item(68, 14).Value = -51665.668

# row 71
$ws.Cells.Item(71, 8).Value = 50167.668
$ws.Cells.Item(71, 10).Value = 50167.668
$ws.Cells.Item(71, 12).Value = 150503.004
$ws.Cells.Item(71, 14).Value = -157991.004

# row 137
$ws.Cells.Item(137, 8).Value = 41897.5
$ws.Cells.Item(137, 10).Value = 41897.5
$ws.Cells.Item(137, 12).Value = 41897.5
$ws.Cells.Item(137, 14).Value = -52097.5

$ws = $wb.Worksheets.Item("CUL")
# row 3
$ws.Cells.Item(3, 8).Value = 2691.125
$ws.Cells.Item(3, 9).Value = 2389.8572
$ws.Cells.Item(3, 10).Value = 4800
$ws.Cells.Item(3, 11).Value = 7169.571599999999
$ws.Cells.Item(3, 12).Value = 14400
$ws.Cells.Item(3, 13).Value = -7057.571599999999
$ws.Cells.Item(3, 14).Value = -14624

# row 4
$ws.Cells.Item(4, 8).Value = 150250
$ws.Cells.Item(4, 10).Value = 500
$ws.Cells.Item(4, 12).Value = 1500
$ws.Cells.Item(4, 14).Value = -1724

# row 44
$ws.Cells.Item(44, 8).Value = 830.125
$ws.Cells.Item(44, 10).Value = 1177.6666
$ws.Cells.Item(44, 12).Value = 3532.9998
$ws.Cells.Item(44, 14).Value = -4328.9998

# row 64
$ws.Cells.Item(64, 8).Value = 1613.7142
$ws.Cells.Item(64, 9).Value = 896
$ws.Cells.Item(64, 10).Value = 1733.3334
$ws.Cells.Item(64, 11).Value = 2688
$ws.Cells.Item(64, 12).Value = 5200.0002
$ws.Cells.Item(64, 13).Value = -2418
$ws.Cells.Item(64, 14).Value = -5740.0002

# row 67
$ws.Cells.Item(67, 8).Value = 1613.7142
$ws.Cells.Item(67, 9).Value = 896
$ws.Cells.Item(67, 10).Value = 1733.3334
$ws.Cells.Item(67, 11).Value = 2688
$ws.Cells.Item(67, 12).Value = 5200.0002
$ws.Cells.Item(67, 13).Value = -1752
$ws.Cells.Item(67, 14).Value = -7072.0002

# row 92
$ws.Cells.Item(92, 8).Value = 746.625
$ws.Cells.Item(92, 10).Value = 661
$ws.Cells.Item(92, 12).Value = 1983
$ws.Cells.Item(92, 14).Value = -4479

$ws = $wb.Worksheets.Item("GSM")
# row 4
$ws.Cells.Item(4, 8).Value = 30000
$ws.Cells.Item(4, 10).Value = 30000
$ws.Cells.Item(4, 12).Value = 30000
$ws.Cells.Item(4, 14).Value = -30224

# row 43
$ws.Cells.Item(43, 8).Value = 17038.9
$ws.Cells.Item(43, 10).Value = 23805.572
$ws.Cells.Item(43, 12).Value = 23805.572
$ws.Cells.Item(43, 14).Value = -24107.572

# row 46
$ws.Cells.Item(46, 8).Value = 34100.8
$ws.Cells.Item(46, 10).Value = 37626
$ws.Cells.Item(46, 12).Value = 37626
$ws.Cells.Item(46, 14).Value = -37938

# row 57
$ws.Cells.Item(57, 8).Value = 38196.6
$ws.Cells.Item(57, 10).Value = 38196.6
$ws.Cells.Item(57, 12).Value = 38196.6
$ws.Cells.Item(57, 14).Value = -39836.6

# row 124
$ws.Cells.Item(124, 8).Value = 41824
$ws.Cells.Item(124, 10).Value = 41824
$ws.Cells.Item(124, 12).Value = 41824
$ws.Cells.Item(124, 14).Value = -51644

# row 132
$ws.Cells.Item(132, 8).Value = 3658.2368
$ws.Cells.Item(132, 9).Value = 2990.3547
$ws.Cells.Item(132, 10).Value = 6616
$ws.Cells.Item(132, 11).Value = 8971.0641
$ws.Cells.Item(132, 12).Value = 19848
$ws.Cells.Item(132, 13).Value = -6441.0641
$ws.Cells.Item(132, 14).Value = -24908

# row 137
$ws.Cells.Item(137, 8).Value = 40510
$ws.Cells.Item(137, 10).Value = 40510
$ws.Cells.Item(137, 12).Value = 40510
$ws.Cells.Item(137, 14).Value = -50710

$ws = $wb.Worksheets.Item("LTW")
# row 87
$ws.Cells.Item(87, 8).Value = 38000
$ws.Cells.Item(87, 9).Value = 2000
$ws.Cells.Item(87, 11).Value = 2000
$ws.Cells.Item(87, 13).Value = -877

# row 90
$ws.Cells.Item(90, 8).Value = 38000
$ws.Cells.Item(90, 9).Value = 2000
$ws.Cells.Item(90, 11).Value = 6000
$ws.Cells.Item(90, 13).Value = -384

# row 136
$ws.Cells.Item(136, 8).Value = 4146.8096
$ws.Cells.Item(136, 9).Value = 1614.0769
$ws.Cells.Item(136, 11).Value = 4842.2307
$ws.Cells.Item(136, 13).Value = -2292.2307

$ws = $wb.Worksheets.Item("WVR")
# row 4
$ws.Cells.Item(4, 8).Value = 51481
$ws.Cells.Item(4, 9).Value = 67474.664
$ws.Cells.Item(4, 10).Value = 3500
$ws.Cells.Item(4, 11).Value = 67474.664
$ws.Cells.Item(4, 12).Value = 3500
$ws.Cells.Item(4, 13).Value = -67361.664
$ws.Cells.Item(4, 14).Value = -3726

# row 123
$ws.Cells.Item(123, 8).Value = 39950
$ws.Cells.Item(123, 10).Value = 39950
$ws.Cells.Item(123, 12).Value = 39950
$ws.Cells.Item(123, 14).Value = -49750

# row 125
$ws.Cells.Item(125, 8).Value = 40055.5
$ws.Cells.Item(125, 10).Value = 40055.5
$ws.Cells.Item(125, 12).Value = 40055.5
$ws.Cells.Item(125, 14).Value = -49895.5
